# Generate Report for Handoff
#
# The f2b59921-a579-48db-89d5-220d41f24c09.md file (last row, row 7, on each
# sheet) just had a new handoff xliff generated, so its "Latest HO Xliff
# Generate Date" / "Latest Handoff Datetime" timestamps move forward.
#
# Overview sheet, column G = "Latest HO Xliff Generate Date"
# zh-cn / de-de sheets, column H = "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-27 18:49:10"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-27 18:49:06"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-27 18:49:10"
